$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.099.36"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "1.823.88"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.10"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4628"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07306"
$ws.Range("E9").Value = "  -2.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8707"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "1.872.85"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.350"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.39"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.492"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008630"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.009"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "27.413.93"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "2.095.52"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.11"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.870"
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("E27").Value = "  -2.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.092"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.32"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.087"
$ws.Range("E30").Value = "  -4.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08906"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.959"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7335"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.458"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.137"
$ws.Range("E35").Value = "  -3.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.476"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05250"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01916"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.926"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.135"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5203"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1630"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.272"
$ws.Range("E45").Value = "  -3.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4870"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("E48").Value = "  -3.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.73"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -3.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06258"
$ws.Range("E51").Value = "  -1.07%  "
